# Se arregla el exportar informes tanto en txt y en excel, se mojoran los informes.
#
# This script fixes several "test"/erroneous values left in the error-report
# workbook so the exported txt/xlsx reports are correct:
#   - A2:A4  "RES00" (placeholder code) -> "SAP035" (real code)
#   - D2     "NONE" (missing first name) -> " jose" (actual first name)
#   - C2     wrong population-group code 99 -> 13
#   - J2     wrong birth date serial 28446 -> 31693
#   - AD2    TFG value out of the flagged range 300 -> 301 (now flagged)
#   - AV2    stray text "r" (flagged) -> correct numeric value 1 (unflagged)
#   - CQ2    wrong professional code 9 -> 1
#   - M2/M3/M4 address cells no longer flagged as erroneous (highlight removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string level content fixes -------------------------------------

# "RES00" -> "SAP035" (affects A2, A3, A4, which all shared this string)
$ws.Range("A2:A4").Value = "SAP035"

# First name for the row-2 patient was a placeholder "NONE"; it's now "jose"
$ws.Range("D2").Value = " jose"

# --- Plain value fixes -------------------------------------------------------

$ws.Range("C2").Value = 13
$ws.Range("J2").Value = 31693
$ws.Range("AD2").Value = 301
$ws.Range("CQ2").Value = 1

# AV2 held a stray string "r" flagged in yellow; replace with the correct
# numeric value (matches AV3/AV4) and drop the error highlight.
$ws.Range("AV2").Value = 1
$ws.Range("AV2").ClearFormats()

# AD2 is now the one flagged as an error (value out of expected range).
$ws.Range("AD2").Interior.Color = 65535

# M2/M3/M4 (addresses) were previously flagged in yellow; the data was
# corrected upstream so the highlight is cleared here.
$ws.Range("M2").ClearFormats()
$ws.Range("M3").ClearFormats()
$ws.Range("M4").ClearFormats()
